$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 2-40: refreshed Price (D) / Volume(1h) (E) quote values.
# D-column cells are stored as text in the source sheet (e.g. thousand-dot
# separators, trailing zeros) so force text format before assigning to avoid
# Excel auto-converting the string into a number.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '29.437.38'
$ws.Range("E2").Value = '  -0.16%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.916.35'
$ws.Range("E3").Value = '  +0.65%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.013'
$ws.Range("E4").Value = '  +0.77%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '325.32'
$ws.Range("E5").Value = '  -0.14%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.011'
$ws.Range("E6").Value = '  +0.63%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4805'
$ws.Range("E7").Value = '  -0.76%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.4046'
$ws.Range("E8").Value = '  -0.39%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.08200'
$ws.Range("E9").Value = '  +0.98%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.007'
$ws.Range("E10").Value = '  -0.07%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '23.31'
$ws.Range("E11").Value = '  -0.78%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.914.20'
$ws.Range("E12").Value = '  +0.49%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '6.046'
$ws.Range("E13").Value = '  +0.48%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '7.222'
$ws.Range("E14").Value = '  +1.95%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '91.11'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.06858'
$ws.Range("E16").Value = '  +1.81%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.011'
$ws.Range("E17").Value = '  +0.52%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.00001038'
$ws.Range("E18").Value = '  -0.07%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '17.51'
$ws.Range("E19").Value = '  -0.87%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '1.010'
$ws.Range("E20").Value = '  +0.57%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '29.453.31'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.659'
$ws.Range("E22").Value = '  +1.68%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '11.82'
$ws.Range("E23").Value = '  +0.24%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.192'
$ws.Range("E24").Value = '  +1.44%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.140.83'
$ws.Range("E25").Value = '  +0.00%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '6.597'
$ws.Range("E26").Value = '  +6.11%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '155.66'
$ws.Range("E27").Value = '  +1.27%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '19.93'
$ws.Range("E28").Value = '  -0.61%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.096'
$ws.Range("E29").Value = '  -0.14%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '120.32'
$ws.Range("E30").Value = '  +1.17%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.012'
$ws.Range("E31").Value = '  -2.08%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.09601'
$ws.Range("E32").Value = '  +0.60%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.595'
$ws.Range("E33").Value = '  +1.41%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.554'
$ws.Range("E34").Value = '  +0.03%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.368'
$ws.Range("E35").Value = '  -1.56%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.06285'
$ws.Range("E36").Value = '  +3.12%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.02281'
$ws.Range("E37").Value = '  +0.83%  '
$ws.Range("E38").Value = '  +0.61%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.5922'
$ws.Range("E39").Value = '  +0.02%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '10.72'
$ws.Range("E40").Value = '  +3.88%  '

# Rows 41-51: a new "Frax" entry was inserted at row 41, shifting the
# previously-listed coins (FraxShare..Quant..MXToken) down one row each;
# the former last row (51, Aave) drops off the bottom of the 51-row table.
# Rewrite B/C/D/E for rows 41-51 directly to their final values, forcing
# text format on column D so numeric-looking strings keep their exact text.
$ws.Range("B41").Value = 'Frax'
$ws.Range("C41").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.010'
$ws.Range("E41").Value = '  +0.49%  '
$ws.Range("B42").Value = 'FraxShare'
$ws.Range("C42").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '7.869'
$ws.Range("E42").Value = '  -0.40%  '
$ws.Range("B43").Value = 'Algorand'
$ws.Range("C43").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.1843'
$ws.Range("E43").Value = '  -0.50%  '
$ws.Range("B44").Value = 'WEMIXToken'
$ws.Range("C44").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.281'
$ws.Range("E44").Value = '  -0.18%  '
$ws.Range("B45").Value = 'RenderToken'
$ws.Range("C45").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.390'
$ws.Range("E45").Value = '  -1.11%  '
$ws.Range("B46").Value = 'EnergySwap'
$ws.Range("C46").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '12.35'
$ws.Range("E46").Value = '  -0.02%  '
$ws.Range("B47").Value = 'Cronos'
$ws.Range("C47").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.07466'
$ws.Range("E47").Value = '  -3.15%  '
$ws.Range("B48").Value = 'Decentraland'
$ws.Range("C48").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.5548'
$ws.Range("E48").Value = '  +0.01%  '
$ws.Range("B49").Value = 'NEARProtocol'
$ws.Range("C49").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.927'
$ws.Range("E49").Value = '  -0.63%  '
$ws.Range("B50").Value = 'Quant'
$ws.Range("C50").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '117.80'
$ws.Range("E50").Value = '  +2.41%  '
$ws.Range("B51").Value = 'MXToken'
$ws.Range("C51").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.429'
$ws.Range("E51").Value = '  +3.71%  '
